$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos price/volume refresh (scheduled GitHub Actions data update).
# Each line below is "CellRef||NewValue". Columns D (Price) and E
# (Volume(1h)) hold values that look numeric/percentage (e.g. "30.233.46",
# "1.000", "  -3.11%  ") but must stay plain text exactly as scraped from
# the source site, so those are written with a leading apostrophe to stop
# Excel's automatic number/date conversion. Columns B/C (Coin name / Link)
# are ordinary text and are written as-is. Rows 19/20 and 27/28 also swap
# their Coin/Link text because the underlying ranking order changed.
$data = @"
D2||30.233.46
E2||  -3.11%  
D3||1.927.38
E3||  -2.94%  
D4||1.000
E4||  +0.09%  
D5||246.32
E5||  -2.79%  
D6||0.7182
E6||  -11.06%  
D7||0.9995
E7||  +0.15%  
D8||0.3252
E8||  -5.64%  
D9||26.44
E9||  +3.18%  
D10||0.06816
E10||  -2.02%  
D11||0.8020
E11||  -4.67%  
D12||0.07935
E12||  -2.15%  
D13||1.927.68
E13||  -2.91%  
D14||5.394
E14||  -1.96%  
D15||94.32
E15||  -6.12%  
E16||  +4.09%  
D17||260.51
E17||  -3.99%  
D18||30.245.82
E18||  -3.10%  
B19||Uniswap
C19||https://coinranking.com/coin/_H5FVG9iW+uniswap-uni
D19||5.833
E19||  +0.61%  
B20||ShibaInu
C20||https://coinranking.com/coin/xz24e0BjL+shibainu-shib
D20||0.000007948
E20||  -0.06%  
D21||2.180.29
E21||  -2.67%  
D22||0.9995
E22||  +0.14%  
D23||0.9994
E23||  +0.03%  
D24||6.857
E24||  -0.98%  
D25||9.662
E25||  -0.62%  
D26||160.17
E26||  -2.49%  
B27||EthereumClassic
C27||https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc
D27||18.92
E27||  -5.17%  
B28||Stellar
C28||https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm
D28||0.1329
E28||  -10.87%  
D29||2.289
E29||  +4.71%  
E30||  +0.68%  
D31||1.548
E31||  -1.21%  
D32||4.418
E32||  -3.26%  
D33||4.192
E33||  -2.67%  
D34||0.05069
E34||  -1.73%  
D35||1.196
E35||  -1.51%  
D36||0.7404
E36||  -2.22%  
D37||2.725
D38||0.01934
E38||  -3.19%  
D39||2.806
E39||  -3.55%  
D40||79.79
E40||  +2.24%  
D41||6.561
E41||  -0.40%  
D42||0.4452
E42||  -4.87%  
D43||2.001
E43||  -3.11%  
D44||0.9993
E44||  +0.15%  
D45||0.8317
E45||  -2.43%  
D46||102.66
E46||  -1.74%  
D47||9.740
E47||  -2.14%  
D48||7.275
E48||  -2.97%  
D49||36.19
E49||  -1.42%  
D50||1.485
E50||  +2.46%  
D51||0.4103
E51||  -4.42%  
"@

$rows = $data -split "`n" | Where-Object { $_.Trim().Length -gt 0 }
foreach ($row in $rows) {
    $parts = $row -split '\|\|', 2
    $ref = $parts[0].Trim()
    $value = $parts[1]
    $col = $ref.Substring(0, 1)
    if ($col -eq 'D' -or $col -eq 'E') {
        $ws.Range($ref).Value = "'" + $value
    } else {
        $ws.Range($ref).Value = $value
    }
}
